$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Force the cell to be treated as text even though the value looks
    # like a plain number (e.g. "214.53"), then restore the cell's
    # style so no stray formatting is left behind.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "25.901.72"
$ws.Range("E2").Value = "  +0.08%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.636.71"
$ws.Range("E3").Value = "  +0.03%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.27%  "

# Row 5 - BNB
Set-TextValue "D5" "214.53"
$ws.Range("E5").Value = "  -0.23%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.98%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.20%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.0635"
$ws.Range("E9").Value = "  +0.47%  "

# Row 10 - Solana
Set-TextValue "D10" "19.60"
$ws.Range("E10").Value = "  -0.40%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.42%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.862.69"

# Row 13 - Polkadot
$ws.Range("E13").Value = "  -0.48%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.632.09"
$ws.Range("E14").Value = "  +0.22%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.543"
$ws.Range("E15").Value = "  -1.68%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -0.44%  "

# Row 17 - Litecoin
Set-TextValue "D17" "62.60"
$ws.Range("E17").Value = "  -0.50%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "25.925.66"
$ws.Range("E18").Value = "  +0.23%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "193.71"
$ws.Range("E20").Value = "  +1.08%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -1.15%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  -0.58%  "

# Row 23 - Chainlink
Set-TextValue "D23" "6.28"
$ws.Range("E23").Value = "  -0.76%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +0.40%  "

# Row 25 - Monero
Set-TextValue "D25" "143.77"
$ws.Range("E25").Value = "  +0.90%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.17%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +2.56%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  -0.30%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  -0.59%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.27%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +1.20%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  -1.22%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  -0.88%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  -2.72%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  +1.37%  "

# Row 36 - ARBITRUM
$ws.Range("E36").Value = "  -0.65%  "

# Row 37 - Maker
$ws.Range("D37").Value = "1.138.62"
$ws.Range("E37").Value = "  -0.84%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  +0.04%  "

# Row 39 - MXToken
$ws.Range("E39").Value = "  -1.13%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -0.01%  "

# Row 42 - Quant
Set-TextValue "D42" "99.40"
$ws.Range("E42").Value = "  -1.18%  "

# Rows 43/44: TrustWalletToken and FraxShare swap rank order
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D43" "5.42"
$ws.Range("E43").Value = "  -3.72%  "

$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D44" "0.796"
$ws.Range("E44").Value = "  -0.73%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.771.70"
$ws.Range("E45").Value = "  +0.01%  "

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  +3.25%  "

# Row 47 - Aave
Set-TextValue "D47" "56.37"
$ws.Range("E47").Value = "  +1.32%  "

# Row 48 - Cronos
$ws.Range("E48").Value = "  +3.31%  "

# Row 49 - RenderToken
$ws.Range("E49").Value = "  -1.29%  "

# Row 50 - Mantle
$ws.Range("E50").Value = "  -0.35%  "

# Row 51 - EnergySwap
Set-TextValue "D51" "7.66"
$ws.Range("E51").Value = "  +0.71%  "
